$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 9 (shifts existing rows 9+ down to 11+)
$ws.Rows("9:10").Insert()

# Row 7: input/configuration_fxe (was input_efficiency/constant_fxe)
$ws.Range("C7").Value = "input"
$ws.Range("D7").Value = "configuration_fxe"

# Row 8: output/configuration_fxe, value 1 (was output_efficiency/constant_fxe, 0.91)
$ws.Range("C8").Value = "output"
$ws.Range("D8").Value = "configuration_fxe"
$ws.Range("G8").Value = 1

# New row 9: input_efficiency / constant_fxe / elecsupply / 1
$ws.Range("A9").Value = "CHE"
$ws.Range("B9").Value = "conv_transmission_elec"
$ws.Range("C9").Value = "input_efficiency"
$ws.Range("D9").Value = "constant_fxe"
$ws.Range("F9").Value = "elecsupply"
$ws.Range("G9").Value = 1

# New row 10: output_efficiency / constant_fxe / elecdelivered / 0.91
$ws.Range("A10").Value = "CHE"
$ws.Range("B10").Value = "conv_transmission_elec"
$ws.Range("C10").Value = "output_efficiency"
$ws.Range("D10").Value = "constant_fxe"
$ws.Range("F10").Value = "elecdelivered"
$ws.Range("G10").Value = 0.91

$ws.Range("G9").Select()
